# Pilot draft ready - to 1st share
#
# Updates the indicator legend descriptions (column E, "indicator_1") on the
# single worksheet to the revised wording/terminology, and moves the active
# cell selection to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value  = "District Population"
$ws.Range("E7").Value  = "Expected damage to built-up assets from river floods (hectares)"
$ws.Range("E9").Value  = "Expected exposure of agricultural land to river floods (hectares)"
$ws.Range("E10").Value = "Expected exposure of agricultural land to river floods (% of ADM agricultural land)"
$ws.Range("E13").Value = "Expected damage to built-up assets from coastal floods (hectares)"
$ws.Range("E14").Value = "Expected damage to built-up assets from coastal floods (% of ADM built-up area)"
$ws.Range("E15").Value = "Population exposed to medium or high landslide hazard (population count)"
$ws.Range("E16").Value = "Population exposed to medium or high landslide hazard (% of ADM population)"
$ws.Range("E17").Value = "Built-up assets exposed to medium or high landslide hazard (Hectares)"
$ws.Range("E18").Value = "Built-up assets exposed to medium or high landslide hazard (% of ADM built-up area)"
$ws.Range("E19").Value = "Frequency of agricultural drought stress affecting at least 30% of arable land during Season 1/Kharif (percentage of historical period 1984-2022)"
$ws.Range("E20").Value = "Frequency of agricultural drought stress affecting at least 30% of arable land during Season 2/Rabi (percentage of historical period 1984-2022)"
$ws.Range("E21").Value = "Expected exposure to heat stress (population count)"
$ws.Range("E22").Value = "Expected exposure to heat stress (% of ADM population)"
$ws.Range("E23").Value = "Expected increase of mortality from air pollution (population count)"
$ws.Range("E24").Value = "Expected increase of mortality from air pollution (% of ADM population)"

# Matches the saved cursor position recorded in the workbook (activeCell E3)
$ws.Range("E3").Select()
